$wb = $excel.ActiveWorkbook

# -- Worksheets involved --
$wsLibrary  = $wb.Worksheets.Item("Library")
$wsFormula  = $wb.Worksheets.Item("Library_Formula")

# -- "Library_Formula" sheet: new row 17 (new formula "attributeExport") --
# Note: write the new shared-string text "attributeExport" before "ExportUtils"
# so the shared-strings table is built in the same order as the target file.
# (Column D / "Description" is intentionally left blank for this row, so it
# is not touched, to avoid materialising an empty-but-styled cell.)
$rowCells = "A17","B17","C17","E17","F17"
foreach ($addr in $rowCells) {
    $wsFormula.Range($addr).Font.Name = "Trebuchet MS"
    $wsFormula.Range($addr).Font.Size = 10
}

$wsFormula.Range("A17").Value = "CREATE/MODIFY"
$wsFormula.Range("C17").Value = "attributeExport"
$wsFormula.Range("B17").Value = "ExportUtils"
$wsFormula.Range("E17").Value = "String"
$wsFormula.Range("F17").Value = "String"

# -- "Library" sheet: new row 7 (new library "ExportUtils") --
# Column A:C on this sheet carries a column-level style, which Excel would
# otherwise auto-apply to a freshly written cell; A7 must stay on the
# worksheet's default (unstyled) look, same as the other Action cells below
# the header that were typed without the column style (rows 4-6), so reset
# it back to "Normal" after writing the value.
$wsLibrary.Range("A7").Value = "CREATE/MODIFY"
$wsLibrary.Range("A7").Style = "Normal"
$wsLibrary.Range("B7").Value = "ExportUtils"
$wsLibrary.Range("B7").Font.Name = "Trebuchet MS"
$wsLibrary.Range("B7").Font.Size = 10

# -- Update selections on each sheet, then make "Library" the active tab --
$wsFormula.Range("B17").Select()
$wsLibrary.Range("A6:A7").Select()
$wsLibrary.Activate()
